# From 1.1 to 1.2 version
# Rotate the content of the TC1/TC2/TC3/TC4 test-case blocks:
#   TC1 (Periodos Avaliativos)        -> TC3
#   TC2 (Niveis das Competencias)     -> TC1
#   TC3 (Avaliacoes)                  -> TC4
#   TC4 (Competencias (portfolio))    -> TC2
#   TC5 (Perfis de Competencias)      -> unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- TC1 block (rows 10-11) now holds "Niveis das Competencias" content ----
$ws.Range("B10").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B11").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"

# ---- TC2 block (rows 19-20) now holds "Competencias (portfolio)" content ----
$ws.Range("B19").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Competencias (portfolio) a partir do menu inicial"
$ws.Range("D19").Value = "SYSTEM exibe a listagem das Competencias (portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B20").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Competencias (portfolio)"

# ---- TC3 block (rows 28-29) now holds "Periodos Avaliativos" content ----
$ws.Range("B28").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial"
$ws.Range("D28").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B29").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Periodos Avaliativos"

# ---- TC4 block (rows 37-38) now holds "Avaliacoes" content ----
$ws.Range("B37").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Avaliacoes a partir do menu inicial"
$ws.Range("D37").Value = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B38").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Avaliacoes"

# TC5 block (rows 46-47) stays "Perfis de Competencias" - no change needed.
